$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text (coinranking prices use
# "." as a thousands separator, e.g. "3.135.35"), not numbers. Any new
# price string that would otherwise parse as a plain float is entered
# via a temporary Text number format so Excel keeps it as a string,
# then the cell style is reset back to Normal so no stray formatting
# is left behind.

$ws.Cells.Item(2, 4).Value = '68.291.70'
$ws.Cells.Item(2, 5).Value = '  +2.28%  '

$ws.Cells.Item(3, 4).Value = '3.135.35'
$ws.Cells.Item(3, 5).Value = '  +2.00%  '

$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '577.66'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.39%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '180.91'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +6.04%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.08%  '

$ws.Cells.Item(8, 4).Value = '3.135.35'
$ws.Cells.Item(8, 5).Value = '  +2.14%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.521'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +1.95%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.54'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +2.47%  '

$ws.Cells.Item(11, 5).Value = '  +2.07%  '

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.469'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.40%  '

$ws.Cells.Item(13, 5).Value = '  +1.25%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '36.84'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +3.16%  '

$ws.Cells.Item(15, 5).Value = '  +1.03%  '

$ws.Cells.Item(16, 4).Value = '68.186.27'
$ws.Cells.Item(16, 5).Value = '  +2.15%  '

$ws.Cells.Item(17, 4).Value = '3.656.18'
$ws.Cells.Item(17, 5).Value = '  +1.94%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.12'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +1.98%  '

$ws.Cells.Item(19, 4).Value = '3.132.69'
$ws.Cells.Item(19, 5).Value = '  +2.09%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.59'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.34%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '487.10'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.44%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.82'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +1.41%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.697'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +1.42%  '

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '83.90'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.38%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.99'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +2.47%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.62'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +4.47%  '

$ws.Cells.Item(28, 5).Value = '  +0.02%  '

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.15'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +4.23%  '

$ws.Cells.Item(30, 5).Value = '  +4.12%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.64'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.15%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '28.24'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.56%  '

$ws.Cells.Item(33, 5).Value = '  +0.92%  '

$ws.Cells.Item(34, 5).Value = '  +3.89%  '

$ws.Cells.Item(35, 5).Value = '  +0.04%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '48.82'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +2.90%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.66'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.36%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.960'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.30%  '

$ws.Cells.Item(39, 5).Value = '  +7.86%  '

$ws.Cells.Item(40, 5).Value = '  +4.45%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '49.17'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.03%  '

$ws.Cells.Item(42, 5).Value = '  +1.68%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.40'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +1.12%  '

$ws.Cells.Item(44, 2).Value = 'Bittensor'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '399.26'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +8.53%  '

$ws.Cells.Item(45, 2).Value = 'dogwifhat'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.71'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +7.84%  '

$ws.Cells.Item(46, 4).Value = '2.797.69'
$ws.Cells.Item(46, 5).Value = '  +1.43%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.30'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +10.74%  '

$ws.Cells.Item(48, 5).Value = '  +1.31%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '135.53'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.17%  '

$ws.Cells.Item(50, 5).Value = '  +0.02%  '

$ws.Cells.Item(51, 5).Value = '  +9.21%  '
